$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits after "Los datos del
#    catálogo" (it gets relocated later in the document).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Underline several italic section headings (adds <w:u w:val="single"/>
#    to the run/paragraph mark formatting without touching the text).
# ---------------------------------------------------------------------
$headings = @(
    "Exploración de los datos",
    "Obtener conjunto inicial de datos",
    "Limpieza de los datos",
    "Integración de los datos",
    "Formateo de los datos",
    "Seleccionar la técnica del modelo",
    "Generar diseño de prueba"
)

foreach ($h in $headings) {
    $rng = $d.Content
    $rng.Find.Execute($h, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Font.Underline = 1
}

# ---------------------------------------------------------------------
# 3) Insert a comma fix: ". Por el contrario en el " ->
#    ". Por el contrario, en el " -- done as a targeted insertion so the
#    edit produces separate runs, mirroring an in-place edit in Word.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(". Por el contrario en el ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $rng.Start
$matchText = $rng.Text
$offset1 = $matchText.IndexOf("contrario")
$pos1 = $matchStart + $offset1
$pos2 = $pos1 + "contrario".Length

$wordRng = $d.Range($pos1, $pos2)
$wordRng.InsertAfter(",")

$pos2b = $pos2 + 1
$wr2 = $d.Range($pos1, $pos2b)
$wr2.Bold = 1
$wr2.Bold = 0

# ---------------------------------------------------------------------
# 4) Spelling corrections: numero -> número, catalogo -> catálogo (x2).
#    These words were each wrapped in <w:proofErr spellStart/spellEnd>
#    markers that must disappear once the spelling is fixed, while the
#    neighbouring (unrelated) runs must stay untouched/unmerged. We
#    achieve this by replacing a 1-char-wider range (so the zero-width
#    proofErr markers end up strictly inside the deleted span and are
#    dropped), then forcing a run split back at the exact word
#    boundaries via a harmless Bold on/off toggle.
# ---------------------------------------------------------------------
function Fix-Word {
    param($doc, $oldWord, $newWord, $context)

    $rng = $doc.Content
    $rng.Find.Execute($context, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $ctxStart = $rng.Start
    $ctxText = $rng.Text

    $wordOffset = $ctxText.IndexOf($oldWord)
    $wordStart = $ctxStart + $wordOffset
    $wordEnd = $wordStart + $oldWord.Length

    $wide = $doc.Range($wordStart - 1, $wordEnd + 1)
    $wideText = $wide.Text
    $prefix = $wideText.Substring(0, 1)
    $suffix = $wideText.Substring($wideText.Length - 1, 1)
    $wide.Text = $prefix + $newWord + $suffix

    $rng2 = $doc.Content
    $rng2.Find.Execute($context.Replace($oldWord, $newWord), $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $ctx2Start = $rng2.Start
    $ctx2Text = $rng2.Text
    $newOffset = $ctx2Text.IndexOf($newWord)
    $newStart = $ctx2Start + $newOffset
    $newEnd = $newStart + $newWord.Length
    $wr = $doc.Range($newStart, $newEnd)
    $wr.Bold = 1
    $wr.Bold = 0
}

Fix-Word $d "numero" "número" "presenta un numero muy alto"
Fix-Word $d "catalogo" "catálogo" "cruzada del catalogo "
Fix-Word $d "catalogo" "catálogo" "estrellas del catalogo "

# ---------------------------------------------------------------------
# 5) Re-insert the "_GoBack" bookmark inside "Seleccionar la técnica del
#    modelo" heading, right before the word "modelo".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Seleccionar la técnica del modelo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headStart = $rng.Start
$headText = $rng.Text
$offset = $headText.IndexOf("modelo")
$splitPos = $headStart + $offset
$bmRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
